$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows (row 2..10), columns A:T.
# Column layout:
#   A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#   E..T numeric metric columns.
$data = @(
    @{Row=2;  A="ECs";  B="Bmp8a"; C="Bmpr1a"; D="ECs";  E=3; F=1; G=0.3786073333333334; H=1.135822;          I=0.3713290366620658; J=0.3713290366620658; K=3; L=1; M=9.918858999999999;  N=29.756577;          O=0.160764128269069; P=0.160764128269069;  Q=3.755352755699333;  R=33.798174801294;    S=0.05969638887997018; T=0.05969638887997018}
    @{Row=3;  A="ECs";  B="Bmp8a"; C="Bmpr1a"; D="FAPs"; E=3; F=1; G=0.3786073333333334; H=1.135822;          I=0.3713290366620658; J=0.3713290366620658; K=3; L=1; M=29.20351433333333;  N=87.61054300000001;  O=0.47332838627826;   P=0.4733283862782601; Q=11.05666468570511;  R=99.50998217134602;  S=0.1757605737015165;  T=0.1757605737015165}
    @{Row=4;  A="ECs";  B="Bmp8a"; C="Bmpr1a"; D="sCs";  E=3; F=1; G=0.3786073333333334; H=1.135822;          I=0.3713290366620658; J=0.3713290366620658; K=3; L=1; M=22.575837;          N=67.72751099999999;  O=0.3659074854526709; P=0.3659074854526709; Q=8.547377444337998;  R=76.926396999042;    S=0.1358720740805791;  T=0.1358720740805792}
    @{Row=5;  A="FAPs"; B="Bmp8a"; C="Bmpr1a"; D="ECs";  E=3; F=1; G=0.4265683333333333; H=1.279705;          I=0.4183680408212104; J=0.4183680408212104; K=3; L=1; M=9.918858999999999;  N=29.756577;          O=0.160764128269069; P=0.160764128269069;  Q=4.231071152198332;  R=38.079640369785;    S=0.06725857337826018; T=0.06725857337826018}
    @{Row=6;  A="FAPs"; B="Bmp8a"; C="Bmpr1a"; D="FAPs"; E=3; F=1; G=0.4265683333333333; H=1.279705;          I=0.4183680408212104; J=0.4183680408212104; K=3; L=1; M=29.20351433333333;  N=87.61054300000001;  O=0.47332838627826;   P=0.4733283862782601; Q=12.45729443664611;  R=112.115649929815;   S=0.1980254696323007;  T=0.1980254696323007}
    @{Row=7;  A="FAPs"; B="Bmp8a"; C="Bmpr1a"; D="sCs";  E=3; F=1; G=0.4265683333333333; H=1.279705;          I=0.4183680408212104; J=0.4183680408212104; K=3; L=1; M=22.575837;          N=67.72751099999999;  O=0.3659074854526709; P=0.3659074854526709; Q=9.630137162694997;  R=86.67123446425498;  S=0.1530839978106495;  T=0.1530839978106495}
    @{Row=8;  A="sCs";  B="Bmp8a"; C="Bmpr1a"; D="ECs";  E=2; F=0.6666666666666666; G=0.214425;          H=0.6432749999999999; I=0.2103029225167239; J=0.2103029225167238; K=3; L=1; M=9.918858999999999;  N=29.756577;          O=0.160764128269069; P=0.160764128269069;  Q=2.126851341075;     R=19.141662069675;    S=0.03380916601083868; T=0.03380916601083867}
    @{Row=9;  A="sCs";  B="Bmp8a"; C="Bmpr1a"; D="FAPs"; E=2; F=0.6666666666666666; G=0.214425;          H=0.6432749999999999; I=0.2103029225167239; J=0.2103029225167238; K=3; L=1; M=29.20351433333333;  N=87.61054300000001;  O=0.47332838627826;   P=0.4733283862782601; Q=6.261963560925;     R=56.357672048325;    S=0.09954234294444286; T=0.09954234294444286}
    @{Row=10; A="sCs";  B="Bmp8a"; C="Bmpr1a"; D="sCs";  E=2; F=0.6666666666666666; G=0.214425;          H=0.6432749999999999; I=0.2103029225167239; J=0.2103029225167238; K=3; L=1; M=22.575837;          N=67.72751099999999;  O=0.3659074854526709; P=0.3659074854526709; Q=4.840823848724999;  R=43.56741463852499;  S=0.07695141356144231; T=0.07695141356144231}
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Range("A$r").Value = $rowData.A
    $ws.Range("B$r").Value = $rowData.B
    $ws.Range("C$r").Value = $rowData.C
    $ws.Range("D$r").Value = $rowData.D
    $ws.Range("E$r").Value = $rowData.E
    $ws.Range("F$r").Value = $rowData.F
    $ws.Range("G$r").Value = $rowData.G
    $ws.Range("H$r").Value = $rowData.H
    $ws.Range("I$r").Value = $rowData.I
    $ws.Range("J$r").Value = $rowData.J
    $ws.Range("K$r").Value = $rowData.K
    $ws.Range("L$r").Value = $rowData.L
    $ws.Range("M$r").Value = $rowData.M
    $ws.Range("N$r").Value = $rowData.N
    $ws.Range("O$r").Value = $rowData.O
    $ws.Range("P$r").Value = $rowData.P
    $ws.Range("Q$r").Value = $rowData.Q
    $ws.Range("R$r").Value = $rowData.R
    $ws.Range("S$r").Value = $rowData.S
    $ws.Range("T$r").Value = $rowData.T
}
